$wb = $excel.ActiveWorkbook

# --- Update status text: "Ready for handoff" -> "In Translation" ---
# Overview sheet: both the zh-cn (E2) and de-de (F2) status cells hold the
# same status string.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# Per-language detail sheets: "Status" column (column C) for the single
# file row on each sheet.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- Narrow the status columns now that the text is shorter ---
$overview.Columns("E:F").ColumnWidth = 13.4101845877511
$zhcn.Columns("C:C").ColumnWidth = 13.4101845877511
$dede.Columns("C:C").ColumnWidth = 13.4101845877511
